$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "1호선"
$ws2 = $wb.Worksheets.Item(2)   # "2호선"

# ---------------------------------------------------------------------
# Sheet "1호선": row 5 was mislabeled "덕정" - fix it to "교대" and
# update its y-coordinate; add a "연결역" (connecting station) column;
# append two new rows of station data.
# ---------------------------------------------------------------------
$ws1.Range("D1").Value = "연결역"

$ws1.Range("A5").Value = "교대"
$ws1.Range("C5").Value = 170

$ws1.Range("A13").Value = "녹차"
$ws1.Range("B13").Value = 300
$ws1.Range("C13").Value = 420
$ws1.Range("D13").Value = "양주"
$ws1.Range("E13").Value = "가능"

$ws1.Range("A14").Value = "골절"
$ws1.Range("B14").Value = 250
$ws1.Range("C14").Value = 470

# ---------------------------------------------------------------------
# Sheet "2호선": replace the computed C5 formula with its static value.
# ---------------------------------------------------------------------
$ws2.Range("C5").Value = 170

# ---------------------------------------------------------------------
# Restore each sheet's last selection, then switch the active tab back
# to "1호선" (it had drifted to "2호선" before this edit).
# ---------------------------------------------------------------------
$ws2.Range("A12").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("F6").Select() | Out-Null
